# Update "想去人数" (interested-count) figures across sheets to reflect
# newly generated data (gh-pages output refresh).
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 501
$ws.Range("F5").Value  = 8542
$ws.Range("F6").Value  = 11
$ws.Range("F7").Value  = 10705
$ws.Range("F8").Value  = 91
$ws.Range("F13").Value = 117
$ws.Range("F22").Value = 1817
$ws.Range("F24").Value = 556
$ws.Range("F27").Value = 65
$ws.Range("F28").Value = 584
$ws.Range("F30").Value = 1184
$ws.Range("F33").Value = 1415
$ws.Range("F35").Value = 342
$ws.Range("F36").Value = 286
$ws.Range("F39").Value = 513
$ws.Range("F42").Value = 791
$ws.Range("F45").Value = 106
$ws.Range("F46").Value = 100

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value  = 37
$ws.Range("F11").Value = 65
$ws.Range("F16").Value = 49
$ws.Range("F17").Value = 383

# --- Sheet: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 501
$ws.Range("F5").Value  = 37
$ws.Range("F10").Value = 8542
$ws.Range("F11").Value = 11
$ws.Range("F12").Value = 10705
$ws.Range("F13").Value = 91
$ws.Range("F15").Value = 117
$ws.Range("F19").Value = 1817
$ws.Range("F21").Value = 556
$ws.Range("F23").Value = 65
$ws.Range("F25").Value = 584
$ws.Range("F28").Value = 1184
$ws.Range("F33").Value = 65
$ws.Range("F34").Value = 1415
$ws.Range("F37").Value = 342
$ws.Range("F39").Value = 513
$ws.Range("F43").Value = 791
$ws.Range("F45").Value = 49
$ws.Range("F46").Value = 383
$ws.Range("F48").Value = 106
$ws.Range("F49").Value = 100

$wb.Save()
